$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the sheet (tab) name to reflect the new "through" date
$ws.Name = "Through 2022-07-13"

# Update the header label in I1 to reflect the new "through" date
$ws.Range("I1").Value = "2022 (through 07-13)"

# Update the July value (row 8) and the recalculated Total (row 14)
$ws.Range("I8").Value = 71
$ws.Range("I14").Value = 877
